$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 33 (shifts existing rows 33.. down by one,
# Excel will extend the sheet dimension automatically and copy formatting
# from the row above for the newly inserted row).
$ws.Rows("33:33").Insert()

# Populate the newly inserted row 33 with the new weekly data point.
$ws.Range("A33").Value = 6
$ws.Range("B33").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C33").Value = "Metropolitana"
$ws.Range("D33").Value = 44721
$ws.Range("E33").Value = 13
$ws.Range("F33").Value = 100114007
$ws.Range("G33").Value = "Jengibre"
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 190
$ws.Range("K33").Value = 11000
$ws.Range("L33").Value = 12000
$ws.Range("M33").Value = 11368
$ws.Range("N33").Value = "$/caja 13 kilos"
$ws.Range("O33").Value = "Perú"
$ws.Range("P33").Value = 874
$ws.Range("Q33").Value = 13
$ws.Range("R33").Value = "Hortaliza"
